$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70, pushing existing rows 70:88 down to 71:89.
$ws.Rows.Item(70).Insert()

# Copy the style of the Fecha column (D) from the row above so the new
# row's date cell keeps the date number format.
$ws.Range("D69").Copy()
$ws.Range("D70").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row 70 with the inserted record's values.
$ws.Cells.Item(70, 1).Value = 5
$ws.Cells.Item(70, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(70, 3).Value = "Maule"
$ws.Cells.Item(70, 4).Value = 44508
$ws.Cells.Item(70, 5).Value = 7
$ws.Cells.Item(70, 6).Value = "Fruta"
$ws.Cells.Item(70, 7).Value = 100108
$ws.Cells.Item(70, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(70, 9).Value = 100108002
$ws.Cells.Item(70, 10).Value = "Mango"
$ws.Cells.Item(70, 11).Value = "Sin especificar"
$ws.Cells.Item(70, 12).Value = "Primera"
$ws.Cells.Item(70, 13).Value = 210
$ws.Cells.Item(70, 14).Value = 7000
$ws.Cells.Item(70, 15).Value = 7000
$ws.Cells.Item(70, 16).Value = 7000
$ws.Cells.Item(70, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(70, 18).Value = "Perú"
$ws.Cells.Item(70, 19).Value = 1750
$ws.Cells.Item(70, 20).Value = 4
